# Update the "Förändrad" (Changed) date column (C) for rows 2-29
# from serial date 45574 (2024-10-09) to 45575 (2024-10-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45574) {
        $cell.Value2 = 45575
    }
}
